$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "row2_c.PNG"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "row2_d.PNG"
$ws.Range("D5").Value = 14

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "row2_e.PNG"

$ws.Range("A7").Select()
